$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20, Column C ("Run Mode") flips from "YES" to "No"
$ws.Range("C20").Value = "No"

# New row 24: Shipper Contact test case
# (shared-string pool order matches D24, then A24, then B24)
$ws.Range("D24").Value = "Shipper Contact Validated Successfully"
$ws.Range("A24").Value = "ShipperContacts_TC001"
$ws.Range("B24").Value = "Validate Shipper Contact.
a) Launch and login applcation.
b)Open Shipper Contact Panel.
c) Click on add new Shipper Contact button.
d) Fill Mandatory Details and click on Add.
e) Verify Newly added Shipper contact in Grid.
f) Edit the Email id and verify in grid.
g)Inactivate the Shipper Contact."
$ws.Range("C24").Value = "YES"

$ws.Range("A24").VerticalAlignment = -4108
$ws.Range("B24").WrapText = $true
$ws.Range("C24").VerticalAlignment = -4108
$ws.Range("D24").VerticalAlignment = -4108

$ws.Rows.Item(24).RowHeight = 120

# Scroll the view down so row 23 is at the top and select the newly-added cell
$excel.ActiveWindow.ScrollRow = 23
$ws.Range("D24").Select()
